$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.289.26"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").Value = "1.833.00"
$ws.Range("E3").Value = "  -0.35%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.63%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "235.79"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.6007"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.98%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.59%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.07018"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -5.28%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.2779"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.73%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "23.53"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -5.42%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07651"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.75%  "

$ws.Range("D12").Value = "1.831.30"
$ws.Range("E12").Value = "  -0.45%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.778"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.71%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.6259"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -7.05%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.000009635"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -5.61%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "78.89"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.36%  "

$ws.Range("D17").Value = "29.312.71"
$ws.Range("E17").Value = "  -0.18%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "5.763"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -7.31%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "223.43"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -4.17%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.50%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.65"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -5.18%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.973"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.41%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.008"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.56%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "156.60"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.1299"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.24%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.975"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -5.88%  "

$ws.Range("E27").Value = "  -3.99%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.06815"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -4.53%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.465"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.89%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.448"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.06%  "

$ws.Range("E31").Value = "  -4.50%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.772"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -6.56%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.104"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.98%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.724"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.96%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.6423"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -8.61%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.553"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.74%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.753"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.05%  "

$ws.Range("D38").Value = "1.211.19"
$ws.Range("E38").Value = "  -1.80%  "

$ws.Range("E39").Value = "  -4.08%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.506"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -5.83%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.9051"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.28%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.60%  "

$ws.Range("D43").Value = "1.993.83"
$ws.Range("E43").Value = "  -0.62%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "100.46"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "62.49"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -4.14%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00000000114"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.02%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "8.523"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.40%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.578"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -7.27%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.4564"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.21%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "6.379"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -8.09%  "
